# add number 5 in blue
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the new value in E1 (extends the row from A1:D1 to A1:E1)
$ws.Range("E1").Value = 5

# Color it blue using the theme's Accent1 color (xlThemeColorAccent1 = 5)
$ws.Range("E1").Font.ThemeColor = 5

# Match the active selection on the sheet to the newly added cell
$ws.Range("E1").Select() | Out-Null
